# The sheet holds 4 years (2014-2017) of monthly data in rows 2..49,
# columns A (label "YYYY-MM") .. D (values), 12 rows per year in
# Jan..Dec order. This edit re-sorts each year's 12-row block so that
# Oct, Nov, Dec come first, followed by Jan..Sep (i.e. a 9-row
# left-rotation of each 12-row block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$rowsPerYear = 12
$yearCount = 4
$lastCol = "D"

for ($y = 0; $y -lt $yearCount; $y++) {
    $blockStart = $firstDataRow + ($y * $rowsPerYear)
    $blockEnd = $blockStart + $rowsPerYear - 1

    # Read the whole 12-row x 4-col block into memory.
    $block = $ws.Range("A" + $blockStart + ":" + $lastCol + $blockEnd).Value2

    # Rows 10..12 of the block (Oct, Nov, Dec) move to the front,
    # rows 1..9 (Jan..Sep) follow, in the same relative order.
    for ($i = 0; $i -lt $rowsPerYear; $i++) {
        $srcRow = (($i + 9) % $rowsPerYear) + 1
        $destRow = $blockStart + $i

        $ws.Cells.Item($destRow, 1).Value = $block[$srcRow, 1]
        $ws.Cells.Item($destRow, 2).Value = $block[$srcRow, 2]
        $ws.Cells.Item($destRow, 3).Value = $block[$srcRow, 3]
        $ws.Cells.Item($destRow, 4).Value = $block[$srcRow, 4]
    }
}
